# Update the "Data" worksheet values per the commit:
#   D6:D10, D16:D17   60000257 -> 60000273
#   E11:E15, E18:E19  202      -> 216

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Column D rows that hold the old id 60000257
$dRows = @(6, 7, 8, 9, 10, 16, 17)
foreach ($r in $dRows) {
    $ws.Cells.Item($r, 4).Value = 60000273
}

# Column E rows that hold the old value 202
$eRows = @(11, 12, 13, 14, 15, 18, 19)
foreach ($r in $eRows) {
    $ws.Cells.Item($r, 5).Value = 216
}
